# Auto-generated edit script: applies scheduled-runner value updates to Sheets/Shinryu_Profits.xlsx
# (workbook tabs: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 226.54546
$ws.Range("I2").Value = 210
$ws.Range("K2").Value = 210
$ws.Range("M2").Value = -97

$ws.Range("H55").Value = 3636448
$ws.Range("I55").Value = 4444525.5
$ws.Range("J55").Value = 100
$ws.Range("K55").Value = 4444525.5
$ws.Range("L55").Value = 100
$ws.Range("M55").Value = -4444311.5
$ws.Range("N55").Value = -528

$ws.Range("H64").Value = 3899.532
$ws.Range("I64").Value = 3773.5652
$ws.Range("J64").Value = 4020.25
$ws.Range("K64").Value = 3773.5652
$ws.Range("L64").Value = 4020.25
$ws.Range("M64").Value = -3525.5652
$ws.Range("N64").Value = -4516.25

$ws.Range("H67").Value = 3899.532
$ws.Range("I67").Value = 3773.5652
$ws.Range("J67").Value = 4020.25
$ws.Range("K67").Value = 3773.5652
$ws.Range("L67").Value = 4020.25
$ws.Range("M67").Value = -2915.5652
$ws.Range("N67").Value = -5736.25

$ws.Range("H132").Value = 2451.125
$ws.Range("I132").Value = 1752.6389
$ws.Range("J132").Value = 8737.5
$ws.Range("K132").Value = 5257.9167
$ws.Range("L132").Value = 26212.5
$ws.Range("M132").Value = -2727.9167
$ws.Range("N132").Value = -31272.5

$ws.Range("H137").Value = 11067128
$ws.Range("I137").Value = 18180272
$ws.Range("J137").Value = 2237.9443
$ws.Range("K137").Value = 54540816
$ws.Range("L137").Value = 6713.8329
$ws.Range("M137").Value = -54538266
$ws.Range("N137").Value = -11813.8329

$ws.Range("H138").Value = 2758.91
$ws.Range("I138").Value = 948.9375
$ws.Range("J138").Value = 3103.6667
$ws.Range("K138").Value = 2846.8125
$ws.Range("L138").Value = 9311.000100000001
$ws.Range("M138").Value = 2293.1875
$ws.Range("N138").Value = -19591.0001

$ws.Range("H141").Value = 1049.6428
$ws.Range("I141").Value = 958.5
$ws.Range("J141").Value = 1277.5
$ws.Range("K141").Value = 2875.5
$ws.Range("L141").Value = 3832.5
$ws.Range("M141").Value = 2304.5
$ws.Range("N141").Value = -14192.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9354.969999999999
$ws.Range("I32").Value = 7088.6914
$ws.Range("J32").Value = 19016.475
$ws.Range("K32").Value = 7088.6914
$ws.Range("L32").Value = 19016.475
$ws.Range("M32").Value = -6801.6914
$ws.Range("N32").Value = -19590.475

$ws.Range("H132").Value = 1204
$ws.Range("I132").Value = 739.5952
$ws.Range("J132").Value = 3642.125
$ws.Range("K132").Value = 2218.7856
$ws.Range("L132").Value = 10926.375
$ws.Range("M132").Value = 311.2143999999998
$ws.Range("N132").Value = -15986.375

$ws.Range("H135").Value = 37407.125
$ws.Range("J135").Value = 37407.125
$ws.Range("L135").Value = 37407.125
$ws.Range("N135").Value = -47547.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1606.7142
$ws.Range("I105").Value = 1606.7142
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1606.7142
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 140.2858000000001
$ws.Range("N105").ClearContents()

$ws.Range("H134").Value = 2427.3333
$ws.Range("I134").Value = 2523.3157
$ws.Range("J134").Value = 2062.6
$ws.Range("K134").Value = 7569.9471
$ws.Range("L134").Value = 6187.799999999999
$ws.Range("M134").Value = -5034.9471
$ws.Range("N134").Value = -11257.8

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 670083.4
$ws.Range("I31").Value = 2067.2917
$ws.Range("J31").Value = 1672107.5
$ws.Range("K31").Value = 2067.2917
$ws.Range("L31").Value = 1672107.5
$ws.Range("M31").Value = -1772.2917
$ws.Range("N31").Value = -1672697.5

$ws.Range("H34").Value = 670083.4
$ws.Range("I34").Value = 2067.2917
$ws.Range("J34").Value = 1672107.5
$ws.Range("K34").Value = 2067.2917
$ws.Range("L34").Value = 1672107.5
$ws.Range("M34").Value = -1865.2917
$ws.Range("N34").Value = -1672511.5

$ws.Range("H58").Value = 2207.7097
$ws.Range("I58").Value = 1450.5454
$ws.Range("J58").Value = 4058.5557
$ws.Range("K58").Value = 1450.5454
$ws.Range("L58").Value = 4058.5557
$ws.Range("M58").Value = -1247.5454
$ws.Range("N58").Value = -4464.5557

$ws.Range("H62").Value = 2980
$ws.Range("I62").Value = 2533.3333
$ws.Range("J62").Value = 3650
$ws.Range("K62").Value = 2533.3333
$ws.Range("L62").Value = 3650
$ws.Range("M62").Value = -1909.3333
$ws.Range("N62").Value = -4898

$ws.Range("H65").Value = 2980
$ws.Range("I65").Value = 2533.3333
$ws.Range("J65").Value = 3650
$ws.Range("K65").Value = 12666.6665
$ws.Range("L65").Value = 18250
$ws.Range("M65").Value = -9546.666499999999
$ws.Range("N65").Value = -24490

$ws.Range("H132").Value = 1662.638
$ws.Range("I132").Value = 1253.6666
$ws.Range("J132").Value = 2502.1052
$ws.Range("K132").Value = 3760.9998
$ws.Range("L132").Value = 7506.3156
$ws.Range("M132").Value = -1230.9998
$ws.Range("N132").Value = -12566.3156

$ws.Range("H134").Value = 2276.4285
$ws.Range("I134").Value = 1046.381
$ws.Range("J134").Value = 4121.5
$ws.Range("K134").Value = 3139.143
$ws.Range("L134").Value = 12364.5
$ws.Range("M134").Value = -604.143
$ws.Range("N134").Value = -17434.5

$ws.Range("H136").Value = 2207.7097
$ws.Range("I136").Value = 1450.5454
$ws.Range("J136").Value = 4058.5557
$ws.Range("K136").Value = 4351.6362
$ws.Range("L136").Value = 12175.6671
$ws.Range("M136").Value = -1801.6362
$ws.Range("N136").Value = -17275.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 30.583334
$ws.Range("I2").Value = 35.588234
$ws.Range("J2").Value = 18.428572
$ws.Range("K2").Value = 35.588234
$ws.Range("L2").Value = 18.428572
$ws.Range("M2").Value = 77.411766
$ws.Range("N2").Value = -244.428572

$ws.Range("H6").Value = 20000
$ws.Range("J6").Value = 20000
$ws.Range("L6").Value = 20000
$ws.Range("N6").Value = -20226

$ws.Range("H16").Value = 20000
$ws.Range("J16").Value = 20000
$ws.Range("L16").Value = 20000
$ws.Range("N16").Value = -20500

$ws.Range("H80").Value = 2890.946
$ws.Range("J80").Value = 3274.0667
$ws.Range("L80").Value = 3274.0667
$ws.Range("N80").Value = -5270.066699999999

$ws.Range("H83").Value = 2890.946
$ws.Range("J83").Value = 3274.0667
$ws.Range("L83").Value = 16370.3335
$ws.Range("N83").Value = -26354.3335

$ws.Range("H93").Value = 20000
$ws.Range("J93").Value = 20000
$ws.Range("L93").Value = 20000
$ws.Range("N93").Value = -23744

$ws.Range("H132").Value = 2873.7666
$ws.Range("I132").Value = 2990.3809
$ws.Range("J132").Value = 2601.6667
$ws.Range("K132").Value = 8971.1427
$ws.Range("L132").Value = 7805.000100000001
$ws.Range("M132").Value = -6441.1427
$ws.Range("N132").Value = -12865.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 597.46155
$ws.Range("I22").Value = 478.75
$ws.Range("J22").Value = 699.2143
$ws.Range("K22").Value = 478.75
$ws.Range("L22").Value = 699.2143
$ws.Range("M22").Value = -183.75
$ws.Range("N22").Value = -1289.2143

$ws.Range("H27").Value = 597.46155
$ws.Range("I27").Value = 478.75
$ws.Range("J27").Value = 699.2143
$ws.Range("K27").Value = 478.75
$ws.Range("L27").Value = 699.2143
$ws.Range("M27").Value = -371.75
$ws.Range("N27").Value = -913.2143

$ws.Range("H46").Value = 900
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 900
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 900
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -1276

$ws.Range("H82").Value = 2464.6428
$ws.Range("I82").Value = 2555.7778
$ws.Range("J82").Value = 2300.6
$ws.Range("K82").Value = 2555.7778
$ws.Range("L82").Value = 2300.6
$ws.Range("M82").Value = -2194.7778
$ws.Range("N82").Value = -3022.6

$ws.Range("H85").Value = 2464.6428
$ws.Range("I85").Value = 2555.7778
$ws.Range("J85").Value = 2300.6
$ws.Range("K85").Value = 2555.7778
$ws.Range("L85").Value = 2300.6
$ws.Range("M85").Value = -1307.7778
$ws.Range("N85").Value = -4796.6

$ws.Range("H132").Value = 1700.69
$ws.Range("I132").Value = 1700.69
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5102.07
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2572.07
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 2669.3225
$ws.Range("I136").Value = 1264.4
$ws.Range("J136").Value = 3338.3333
$ws.Range("K136").Value = 3793.2
$ws.Range("L136").Value = 10014.9999
$ws.Range("M136").Value = -1243.2
$ws.Range("N136").Value = -15114.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 915.7143
$ws.Range("I100").Value = 637.5
$ws.Range("J100").Value = 1286.6666
$ws.Range("K100").Value = 1275
$ws.Range("L100").Value = 2573.3332
$ws.Range("M100").Value = -734
$ws.Range("N100").Value = -3655.3332

$ws.Range("H132").Value = 1660.836
$ws.Range("I132").Value = 1102.0392
$ws.Range("J132").Value = 4510.7
$ws.Range("K132").Value = 3306.1176
$ws.Range("L132").Value = 13532.1
$ws.Range("M132").Value = -776.1175999999996
$ws.Range("N132").Value = -18592.1

$ws.Range("H136").Value = 4080.7908
$ws.Range("I136").Value = 4618.4814
$ws.Range("J136").Value = 3173.4375
$ws.Range("K136").Value = 13855.4442
$ws.Range("L136").Value = 9520.3125
$ws.Range("M136").Value = -11305.4442
$ws.Range("N136").Value = -14620.3125
